# Swap the "Recorded By" value order for System + user email:
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# Only cells that contain exactly that combined string are affected;
# cells with just "System" or just the bare email are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
